$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$result = $find.Execute("Berechnung der F", $true, $false, $false, $false, $false, $true, 1, $false, "Berechnungsgrundlagen der f", 2)
Write-Host "Replace1 result: $result"
